$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Travel / Domestic (row 42): supplement the original $3220 with an additional
# $1220 of travel costs for each of the three budget years, expressed as the
# same formula Excel would record ("=3220+1220") so the underlying number is
# still traceable to the original estimate.
$ws.Range("B42").Formula = "=3220+1220"
$ws.Range("C42").Formula = "=3220+1220"
$ws.Range("D42").Formula = "=3220+1220"

# Materials and Supplies (row 47): Years 2 and 3 drop from $4000 to $2000.
$ws.Range("C47").Value = 2000
$ws.Range("D47").Value = 2000

# Recalculate so every dependent subtotal/grand-total formula downstream
# (rows 44, 55, 71, 73, 75, 77, etc.) picks up the new cached values.
$excel.CalculateFull()

# Reflect the author's final on-screen scroll position / selection.
$ws.Activate()
$ws.Range("D57").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
